{"js": "// The \"Threads, Exceptions, C++ Language Standards\" table cell had its\n// text split across three bold Arial runs: \"Threads, \", \"Exceptions, \"\n// and \"C++ Language Standards\". The edit removes the middle\n// \"Exceptions, \" run/text so the cell reads \"Threads, C++ Language\n// Standards\".\nconst results = context.document.body.search(\"Exceptions, \", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].delete();\n  await context.sync();\n}\n", "ps1": "# The \"Threads, Exceptions, C++ Language Standards\" table cell had its\n# text split across three bold Arial runs: \"Threads, \", \"Exceptions, \"\n# and \"C++ Language Standards\". The edit removes the middle\n# \"Exceptions, \" run/text so the cell reads \"Threads, C++ Language\n# Standards\".\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute(\"Exceptions, \")\n\nif ($found) {\n    $rng.Text = \"\"\n}\n"}
